$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 94, shifting the trailing spacer row and the three
# summary rows (sum [min] / sum [h] / sum [working weeks]) down by one.
$ws.Rows("94:94").Insert()

# The entry that used to be the last data row (93) actually ran 45 minutes
# (0:45) longer than originally recorded.
$ws.Range("E93").Value = 0.57291666666666663

# Populate the newly inserted row 94 with the extra working-hours entry.
$ws.Range("A94").Value = 2014
$ws.Range("B94").Value = 3
$ws.Range("C94").Value = 25
$ws.Range("D94").Value = 0.61111111111111105
$ws.Range("E94").Value = 0.75

# Carry the time-spent formulas down into the new row, matching the pattern
# used by every other data row (minutes worked, then hours worked).
$ws.Range("F94").Formula = "=(E94-D94)*24*60"
$ws.Range("G94").Formula = "=F94/60"

# Match the number formats of the surrounding data rows: hh:mm for the
# start/end times, integer minutes, 2-decimal hours.
$ws.Range("D94:E94").NumberFormat = "hh:mm;@"
$ws.Range("F94").NumberFormat = "0"
$ws.Range("G94").NumberFormat = "0.00"

# The spacer row moved from 94 to 95, so the running SUM needs to cover
# through the new last data row, F95.
$ws.Range("F96").Formula = "=SUM(F2:F95)"

# Restore the viewport/selection the author had at save time.
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E95").Select() | Out-Null
